$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'29.660.63"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +2.69%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.862.61"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +2.03%  "
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'0.9993"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +0.03%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'245.77"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  +2.78%  "
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'  +1.50%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'  +0.03%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.07752"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +1.81%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.3072"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +1.93%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'23.68"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +1.30%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.07782"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  +0.81%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'5.163"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +2.43%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'1.854.24"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +1.66%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'92.37"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +2.69%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.6937"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +3.37%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'6.577"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +2.68%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'29.650.93"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +2.72%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'0.000008368"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +1.19%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'2.113.33"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +1.69%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'242.28"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -0.14%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'12.77"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +1.30%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.9999"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'7.629"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +3.48%  "
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'  +0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.1514"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +2.99%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'8.927"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +2.58%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'159.68"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -0.34%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'18.32"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  +1.00%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'1.544"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  +1.09%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'4.264"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +1.91%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'4.195"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +1.58%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'1.194"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  +0.32%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'0.05109"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +0.25%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'0.7874"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +5.29%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'1.904"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +5.18%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'1.159"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +1.64%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'2.687"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +0.14%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'1.332.03"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +11.26%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.01883"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +2.89%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'2.735"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +2.30%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.9607"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  +5.13%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'5.990"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  +14.68%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'106.56"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -1.56%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.9993"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +0.06%  "
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'  +3.70%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'9.783"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  +3.64%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'2.010.53"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +1.65%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'0.5213"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  +1.24%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'65.01"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  +4.52%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'  +3.75%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'7.006"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  +2.15%  "
$c.Style = "Normal"
